# Actualización automática 2025-10-22 14:30:09
$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" ---
$wsVentasGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsVentasGrupo.Range("M21").Value = 9568.26
$wsVentasGrupo.Range("P21").Value = 478.25

# --- Sheet "VENTA MENSUAL" ---
$wsVentaMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsVentaMensual.Range("F21").Value = 10604.54
$wsVentaMensual.Range("F37").Value = 25645.92
# Column widths stored in the XML are ColumnWidth + 5/6 (Excel's internal
# padding for this font), so back-solve the ColumnWidth that yields the
# target stored width.
$wsVentaMensual.Columns.Item(6).ColumnWidth = 14 - (5/6)
$wsVentaMensual.Columns.Item(4).ColumnWidth = 13 - (5/6)

# --- Sheet "CUMPLIMIENTO MENSUAL" ---
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumplimiento.Range("D8").Value = 478.25
$wsCumplimiento.Range("E8").Value = 1.967743214072016
$wsCumplimiento.Range("F8").Value = 0.9959023937747448
$wsCumplimiento.Range("D12").Value = 23611.7
$wsCumplimiento.Range("E12").Value = -1910.43
$wsCumplimiento.Range("F12").Value = 1.088033096680517
$wsCumplimiento.Range("D14").Value = 25645.92
$wsCumplimiento.Range("E14").Value = 10939.64723718182
$wsCumplimiento.Range("F14").Value = 0.7009846214420891
$wsCumplimiento.Columns.Item(4).ColumnWidth = 13 - (5/6)
